# Insert a new weekly price record for "Femacal de La Calera" / Ajo at
# row 199 (pushing the existing rows 199-246 down to 200-247), then
# populate the new row with the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("199:199").Insert()

$ws.Range("A199").Value2 = 3
$ws.Range("B199").Value2 = "Femacal de La Calera"
$ws.Range("C199").Value2 = "Coquimbo"
$ws.Range("D199").Value2 = 44508
$ws.Range("E199").Value2 = 5
$ws.Range("F199").Value2 = 100112003
$ws.Range("G199").Value2 = "Ajo"
$ws.Range("H199").Value2 = "Chino"
$ws.Range("I199").Value2 = "Primera"
$ws.Range("J199").Value2 = 73
$ws.Range("K199").Value2 = 16000
$ws.Range("L199").Value2 = 16500
$ws.Range("M199").Value2 = 16260
$ws.Range("N199").Value2 = "$/caja 10 kilos"
$ws.Range("O199").Value2 = "China"
$ws.Range("P199").Value2 = 1626
$ws.Range("Q199").Value2 = 10
$ws.Range("R199").Value2 = "Hortaliza"
